$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New user rows to append below existing data (rows 4-6)
$data = @(
    @("ara", "wr", "ara@gmail.com", "hfg", "asdfgh"),
    @("arw", "wer", "wer@gmail.com", "awed", "1111111"),
    @("awr", "wer", "a@gmail.com", "fda", "ararara")
)

$startRow = 4
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $values = $data[$i]
    for ($c = 0; $c -lt $values.Length; $c++) {
        $cell = $ws.Cells.Item($row, $c + 1)
        $value = $values[$c]
        if ($value -match '^-?[0-9]+(\.[0-9]+)?$') {
            # Keep purely-numeric-looking text values (e.g. "1111111") stored
            # as text, not auto-converted to a number, then strip the
            # resulting formatting so the cell keeps the default style.
            $cell.NumberFormat = "@"
            $cell.Value = $value
            $cell.ClearFormats()
        } else {
            $cell.Value = $value
        }
    }
}
